$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes contain values that look numeric (e.g. "1.000",
# "0.9995") but must remain plain text, matching the source data which stores
# them as inline strings. Force text format before assignment, then restore the
# default (Normal) style so no stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "27.135.61"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "1.710.20"
$ws.Range("E3").Value = "  -3.04%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.16%  "
Set-TextValue $ws.Range("D5") "307.34"
$ws.Range("E5").Value = "  -6.31%  "
Set-TextValue $ws.Range("D6") "0.9995"
$ws.Range("E6").Value = "  -0.21%  "
Set-TextValue $ws.Range("D7") "0.4787"
$ws.Range("E7").Value = "  +7.16%  "
Set-TextValue $ws.Range("D8") "0.3436"
$ws.Range("E8").Value = "  -3.03%  "
Set-TextValue $ws.Range("D9") "41.86"
$ws.Range("E9").Value = "  -0.18%  "
Set-TextValue $ws.Range("D10") "0.07293"
$ws.Range("E10").Value = "  -1.49%  "
Set-TextValue $ws.Range("D11") "1.048"
$ws.Range("E11").Value = "  -4.74%  "
Set-TextValue $ws.Range("D12") "1.000"
$ws.Range("E12").Value = "  -0.07%  "
Set-TextValue $ws.Range("D13") "19.86"
$ws.Range("E13").Value = "  -4.86%  "
Set-TextValue $ws.Range("D14") "5.847"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "1.708.27"
$ws.Range("E15").Value = "  -3.15%  "
Set-TextValue $ws.Range("D16") "6.829"
$ws.Range("E16").Value = "  -5.65%  "
Set-TextValue $ws.Range("D17") "89.18"
$ws.Range("E17").Value = "  -3.95%  "
Set-TextValue $ws.Range("D18") "0.00001040"
$ws.Range("E18").Value = "  -1.89%  "
Set-TextValue $ws.Range("D19") "0.06348"
$ws.Range("E19").Value = "  -1.16%  "
Set-TextValue $ws.Range("D20") "0.9994"
$ws.Range("E20").Value = "  -0.13%  "
Set-TextValue $ws.Range("D21") "16.44"
$ws.Range("E21").Value = "  -3.67%  "
Set-TextValue $ws.Range("D22") "5.594"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").Value = "27.175.00"
$ws.Range("E23").Value = "  -2.55%  "
Set-TextValue $ws.Range("D24") "10.85"
$ws.Range("E24").Value = "  -3.44%  "
Set-TextValue $ws.Range("D25") "2.100"
$ws.Range("E25").Value = "  -0.28%  "
Set-TextValue $ws.Range("D26") "154.84"
$ws.Range("E26").Value = "  -3.70%  "
$ws.Range("E27").Value = "  -3.52%  "
$ws.Range("D28").Value = "1.909.22"
$ws.Range("E28").Value = "  -2.94%  "
Set-TextValue $ws.Range("D29") "2.079"
$ws.Range("E29").Value = "  -3.04%  "
Set-TextValue $ws.Range("D30") "119.32"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("E31").Value = "  -8.26%  "
Set-TextValue $ws.Range("D32") "0.09247"
$ws.Range("E32").Value = "  +0.60%  "
Set-TextValue $ws.Range("D33") "3.579"
$ws.Range("E33").Value = "  -2.95%  "
Set-TextValue $ws.Range("D34") "5.306"
$ws.Range("E34").Value = "  -6.11%  "
Set-TextValue $ws.Range("D35") "0.02193"
$ws.Range("E35").Value = "  -3.74%  "
Set-TextValue $ws.Range("D36") "0.05853"
$ws.Range("E36").Value = "  -5.44%  "
Set-TextValue $ws.Range("D37") "11.05"
$ws.Range("E37").Value = "  -6.62%  "
Set-TextValue $ws.Range("D38") "0.1991"
$ws.Range("E38").Value = "  -5.11%  "
Set-TextValue $ws.Range("D39") "4.741"
$ws.Range("E39").Value = "  -4.11%  "
Set-TextValue $ws.Range("D40") "0.9993"
Set-TextValue $ws.Range("D41") "1.402"
$ws.Range("E41").Value = "  +0.60%  "
Set-TextValue $ws.Range("D42") "0.5873"
$ws.Range("E42").Value = "  -6.88%  "
Set-TextValue $ws.Range("D43") "1.108"
$ws.Range("E43").Value = "  -6.25%  "
Set-TextValue $ws.Range("D44") "7.447"
$ws.Range("E44").Value = "  -5.06%  "
Set-TextValue $ws.Range("D45") "12.53"
$ws.Range("E45").Value = "  -4.94%  "
Set-TextValue $ws.Range("D46") "3.556"
$ws.Range("E46").Value = "  -5.07%  "
Set-TextValue $ws.Range("D47") "0.5614"
$ws.Range("E47").Value = "  -4.06%  "
Set-TextValue $ws.Range("D48") "117.60"
$ws.Range("E48").Value = "  -3.80%  "
Set-TextValue $ws.Range("D49") "1.840"
$ws.Range("E49").Value = "  -5.68%  "
Set-TextValue $ws.Range("D50") "0.06624"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("E51").Value = "  -4.29%  "
